$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# ---------------------------------------------------------------------
# 1) Rename the header row: "<name>_old" -> "<name>_FV2410",
#    "<name>_new" -> "<name>_FV2504" ("diff" stays untouched).
# ---------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2) Turn the data range A1:U66 into an Excel Table ("Table1") with an
#    autofilter, without dragging in a header-row dxf override (stash
#    the header row's existing formatting on a scratch cell, strip it
#    so the freshly-created table doesn't need to capture an override,
#    then restore it once the table exists).
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")

$ws.Range("A1").Copy() | Out-Null
$ws.Range("A70").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$ws.Range("A70").Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A70").Clear()

# ---------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
